$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# Add the new row of data below the existing table (row 8)
$ws.Range("A8").Value = "Prezentáció"
$ws.Range("B8").Value = "Bankó Olivér"

# Match formatting of the preceding data row (row 7) which already has the
# centered, unshaded style used for all non-header data rows.
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep the same selection Excel had before editing (cell C8)
$ws.Range("C8").Select()
